$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.508609294891357
$ws.Range("B1").Value = 2.131599903106689
$ws.Range("C1").Value = 2.537649154663086
$ws.Range("D1").Value = 3.030050992965698
$ws.Range("E1").Value = 2.32846474647522
